$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume(1h) values per latest data refresh
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "76.012.02"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.57%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.914.76"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.47%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "198.63"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +4.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "598.16"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.69%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.51%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.910.65"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.429"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +15.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.161"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.89"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.433.33"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "75.885.00"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.36"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.910.82"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.84"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -4.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.79"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "377.84"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.19"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.29"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.058.78"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.25%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.63"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.28%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +4.71%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "501.17"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -5.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.75"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.35%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.59%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.12"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.14%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -6.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "180.81"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.343"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.01"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.11%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0905"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +6.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.20"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.92%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.16"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.08%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.70%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.578"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.53%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.659"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +6.54%  "
